$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H21").Value = 300
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 300
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H38").Value = 366.2
$ws.Range("I38").Value = 7.75
$ws.Range("K38").Value = 23.25
$ws.Range("M38").Value = 348.75
$ws.Range("H86").Value = 2186.75
$ws.Range("I86").Value = 849
$ws.Range("K86").Value = 849
$ws.Range("M86").Value = 274
$ws.Range("H89").Value = 2186.75
$ws.Range("I89").Value = 849
$ws.Range("K89").Value = 4245
$ws.Range("M89").Value = 1371
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("N92").ClearContents()
$ws.Range("H107").Value = 2641.5
$ws.Range("I107").Value = 1837.5
$ws.Range("K107").Value = 1837.5
$ws.Range("M107").Value = 82.5
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("H112").Value = 1883.3334
$ws.Range("I112").Value = 1050
$ws.Range("K112").Value = 3150
$ws.Range("M112").Value = -2042
$ws.Range("H137").Value = 1633.6666
$ws.Range("I137").Value = 1545.75
$ws.Range("K137").Value = 4637.25
$ws.Range("M137").Value = -2087.25

$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 3098.5557
$ws.Range("I2").Value = 2426.8572
$ws.Range("J2").Value = 5449.5
$ws.Range("K2").Value = 2426.8572
$ws.Range("L2").Value = 5449.5
$ws.Range("M2").Value = -2313.8572
$ws.Range("N2").Value = -5675.5
$ws.Range("H32").Value = 8252.429
$ws.Range("I32").Value = 8044.6665
$ws.Range("K32").Value = 8044.6665
$ws.Range("M32").Value = -7757.6665
$ws.Range("H97").Value = 3069.25
$ws.Range("I97").Value = 2753.875
$ws.Range("K97").Value = 2753.875
$ws.Range("M97").Value = -2257.875
$ws.Range("H116").Value = 3098.5557
$ws.Range("I116").Value = 2426.8572
$ws.Range("J116").Value = 5449.5
$ws.Range("K116").Value = 2426.8572
$ws.Range("L116").Value = 5449.5
$ws.Range("M116").Value = -132.8571999999999
$ws.Range("N116").Value = -10037.5
$ws.Range("H122").Value = 9814.546
$ws.Range("I122").Value = 9814.546
$ws.Range("K122").Value = 29443.638
$ws.Range("M122").Value = -26993.638
$ws.Range("H132").Value = 791.6667
$ws.Range("I132").Value = 791.6667
$ws.Range("K132").Value = 2375.0001
$ws.Range("M132").Value = 154.9998999999998

$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 3098.5557
$ws.Range("I3").Value = 2426.8572
$ws.Range("J3").Value = 5449.5
$ws.Range("K3").Value = 2426.8572
$ws.Range("L3").Value = 5449.5
$ws.Range("M3").Value = -2312.8572
$ws.Range("N3").Value = -5677.5
$ws.Range("H22").Value = 2958.3333
$ws.Range("I22").Value = 4250
$ws.Range("K22").Value = 4250
$ws.Range("M22").Value = -4077
$ws.Range("H105").Value = 8000
$ws.Range("J105").Value = 10000
$ws.Range("L105").Value = 10000
$ws.Range("N105").Value = -13494
$ws.Range("H107").Value = 3417.0833
$ws.Range("I107").Value = 3478.6667
$ws.Range("J107").Value = 3232.3333
$ws.Range("K107").Value = 3478.6667
$ws.Range("L107").Value = 3232.3333
$ws.Range("M107").Value = -1558.6667
$ws.Range("N107").Value = -7072.3333
$ws.Range("H134").Value = 4531.7334
$ws.Range("I134").Value = 2713
$ws.Range("J134").Value = 8169.2
$ws.Range("K134").Value = 8139
$ws.Range("L134").Value = 24507.6
$ws.Range("M134").Value = -5604
$ws.Range("N134").Value = -29577.6

$ws = $wb.Worksheets.Item(4)
$ws.Range("H16").Value = 350
$ws.Range("I16").Value = 330
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 330
$ws.Range("L16").Value = 400
$ws.Range("M16").Value = -43
$ws.Range("N16").Value = -974
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = 50
$ws.Range("H99").Value = 2999
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 2999
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 2999
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -5995
$ws.Range("H107").Value = 750
$ws.Range("I107").Value = 750
$ws.Range("K107").Value = 750
$ws.Range("M107").Value = 1170
$ws.Range("H113").Value = 350
$ws.Range("I113").Value = 330
$ws.Range("J113").Value = 400
$ws.Range("K113").Value = 330
$ws.Range("L113").Value = 400
$ws.Range("M113").Value = 1840
$ws.Range("N113").Value = -4740
$ws.Range("H126").Value = 2999
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 2999
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 8997
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -13937

$ws = $wb.Worksheets.Item(5)
$ws.Range("H5").Value = 192.5
$ws.Range("H103").Value = 1659.4
$ws.Range("J103").Value = 2248.5
$ws.Range("L103").Value = 6745.5
$ws.Range("N103").Value = -8503.5
$ws.Range("H113").Value = 888
$ws.Range("I113").Value = 888
$ws.Range("K113").Value = 2664
$ws.Range("M113").Value = -494
$ws.Range("H131").Value = 3492.9
$ws.Range("I131").Value = 2030
$ws.Range("J131").Value = 3655.4443
$ws.Range("K131").Value = 6090
$ws.Range("L131").Value = 10966.3329
$ws.Range("M131").Value = -1050
$ws.Range("N131").Value = -21046.3329
$ws.Range("H135").Value = 192.5

$ws = $wb.Worksheets.Item(6)
$ws.Range("H5").Value = 54.5
$ws.Range("I5").Value = 54.5
$ws.Range("K5").Value = 54.5
$ws.Range("M5").Value = 57.5
$ws.Range("H43").Value = 10731
$ws.Range("I43").Value = 5039
$ws.Range("K43").Value = 5039
$ws.Range("M43").Value = -4888
$ws.Range("H46").Value = 11087.5
$ws.Range("I46").Value = 2175
$ws.Range("J46").Value = 20000
$ws.Range("K46").Value = 2175
$ws.Range("L46").Value = 20000
$ws.Range("M46").Value = -2019
$ws.Range("N46").Value = -20312
$ws.Range("H57").Value = 27000
$ws.Range("J57").Value = 27000
$ws.Range("L57").Value = 27000
$ws.Range("N57").Value = -28640
$ws.Range("H113").Value = 6333.3335
$ws.Range("J113").Value = 7000
$ws.Range("L113").Value = 7000
$ws.Range("N113").Value = -11340
$ws.Range("H126").Value = 4124.1665
$ws.Range("I126").Value = 4124.1665
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 12372.4995
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -9902.499500000002
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 2713.6667
$ws.Range("I132").Value = 2176.5715
$ws.Range("K132").Value = 6529.7145
$ws.Range("M132").Value = -3999.7145

$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 3635.4285
$ws.Range("I7").Value = 3408
$ws.Range("K7").Value = 3408
$ws.Range("M7").Value = -3296
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H40").Value = 4557.2
$ws.Range("I40").Value = 4557.2
$ws.Range("K40").Value = 4557.2
$ws.Range("M40").Value = -4421.2
$ws.Range("H46").Value = 4201.6
$ws.Range("J46").Value = 5002
$ws.Range("L46").Value = 5002
$ws.Range("N46").Value = -5378
$ws.Range("H61").Value = 4374.75
$ws.Range("I61").Value = 4333
$ws.Range("K61").Value = 4333
$ws.Range("M61").Value = -4131
$ws.Range("H113").Value = 4374.75
$ws.Range("I113").Value = 4333
$ws.Range("K113").Value = 4333
$ws.Range("M113").Value = -2163
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H126").Value = 3635.4285
$ws.Range("I126").Value = 3408
$ws.Range("K126").Value = 10224
$ws.Range("M126").Value = -7754

$ws = $wb.Worksheets.Item(8)
$ws.Range("H62").Value = 3665.6667
$ws.Range("J62").Value = 4000
$ws.Range("L62").Value = 4000
$ws.Range("N62").Value = -5248
$ws.Range("H65").Value = 3665.6667
$ws.Range("J65").Value = 4000
$ws.Range("L65").Value = 20000
$ws.Range("N65").Value = -26240
$ws.Range("H100").Value = 7882.643
$ws.Range("I100").Value = 9044.637000000001
$ws.Range("K100").Value = 18089.274
$ws.Range("M100").Value = -17548.274
$ws.Range("H113").Value = 599
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 599
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1797
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6137
$ws.Range("H126").Value = 1736.375
$ws.Range("I126").Value = 1761.1428
$ws.Range("J126").Value = 1563
$ws.Range("K126").Value = 5283.428400000001
$ws.Range("L126").Value = 4689
$ws.Range("M126").Value = -2813.428400000001
$ws.Range("N126").Value = -9629
